# Insert a new data row right above the current row 194, pushing the
# existing rows 194:283 down to 195:284 (dimension grows from R283 to R284).
# The new row 194 is a copy of the (now shifted) row below it, with a new
# date and new price/origin figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 194 and below down by one row.
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 - same categorical data as the row
# that used to occupy 194 (now at 195), but with its own date / price /
# origin values.
$ws.Range("A194").Value = 8
$ws.Range("B194").Value = "Terminal La Palmera de La Serena"
$ws.Range("C194").Value = "Coquimbo"
$ws.Range("D194").Value = 44839
$ws.Range("E194").Value = 4
$ws.Range("F194").Value = 100112031
$ws.Range("G194").Value = "Poroto verde"
$ws.Range("H194").Value = "Magnum"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 600
$ws.Range("K194").Value = 34000
$ws.Range("L194").Value = 35000
$ws.Range("M194").Value = 34500
$ws.Range("N194").Value = "`$/malla 25 kilos"
$ws.Range("O194").Value = "Perú"
$ws.Range("P194").Value = 1380
$ws.Range("Q194").Value = 25
$ws.Range("R194").Value = "Hortaliza"
